# Implement slider for slew speed:
# The "Slew Speed" translation row (row 51 on the Translation sheet) is no
# longer needed as a standalone label now that a slider widget is used, so
# remove it and let every row below shift up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Rows.Item(51).Delete()
